$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-10-04 Saturday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-10-05 Sunday", 2) | Out-Null
$d.Content.Find.Execute("216÷3=72, 0", $true, $false, $false, $false, $false, $true, 1, $false, "208÷4=52, 0", 2) | Out-Null
$d.Content.Find.Execute("640÷4=160, 0", $true, $false, $false, $false, $false, $true, 1, $false, "522÷8=65, 2", 2) | Out-Null
$d.Content.Find.Execute("219÷3=73, 0", $true, $false, $false, $false, $false, $true, 1, $false, "736÷2=368, 0", 2) | Out-Null
$d.Content.Find.Execute("593÷6=98, 5", $true, $false, $false, $false, $false, $true, 1, $false, "743÷2=371, 1", 2) | Out-Null
$d.Content.Find.Execute("759÷4=189, 3", $true, $false, $false, $false, $false, $true, 1, $false, "269÷2=134, 1", 2) | Out-Null
$d.Content.Find.Execute("314÷5=62, 4", $true, $false, $false, $false, $false, $true, 1, $false, "570÷7=81, 3", 2) | Out-Null
$d.Content.Find.Execute("227÷3=75, 2", $true, $false, $false, $false, $false, $true, 1, $false, "975÷7=139, 2", 2) | Out-Null
$d.Content.Find.Execute("207÷8=25, 7", $true, $false, $false, $false, $false, $true, 1, $false, "131÷3=43, 2", 2) | Out-Null
$d.Content.Find.Execute("284÷4=71, 0", $true, $false, $false, $false, $false, $true, 1, $false, "244÷9=27, 1", 2) | Out-Null
$d.Content.Find.Execute("477÷4=119, 1", $true, $false, $false, $false, $false, $true, 1, $false, "539÷3=179, 2", 2) | Out-Null
$d.Content.Find.Execute("365÷3=121, 2", $true, $false, $false, $false, $false, $true, 1, $false, "601÷9=66, 7", 2) | Out-Null
$d.Content.Find.Execute("988÷8=123, 4", $true, $false, $false, $false, $false, $true, 1, $false, "409÷7=58, 3", 2) | Out-Null
$d.Content.Find.Execute("329÷9=36, 5", $true, $false, $false, $false, $false, $true, 1, $false, "796÷9=88, 4", 2) | Out-Null
$d.Content.Find.Execute("158÷9=17, 5", $true, $false, $false, $false, $false, $true, 1, $false, "878÷7=125, 3", 2) | Out-Null
$d.Content.Find.Execute("215÷7=30, 5", $true, $false, $false, $false, $false, $true, 1, $false, "226÷4=56, 2", 2) | Out-Null
$d.Content.Find.Execute("241÷4=60, 1", $true, $false, $false, $false, $false, $true, 1, $false, "441÷3=147, 0", 2) | Out-Null
$d.Content.Find.Execute("662÷2=331, 0", $true, $false, $false, $false, $false, $true, 1, $false, "858÷9=95, 3", 2) | Out-Null
$d.Content.Find.Execute("605÷4=151, 1", $true, $false, $false, $false, $false, $true, 1, $false, "540÷2=270, 0", 2) | Out-Null
$d.Content.Find.Execute("143÷7=20, 3", $true, $false, $false, $false, $false, $true, 1, $false, "622÷9=69, 1", 2) | Out-Null
$d.Content.Find.Execute("199÷8=24, 7", $true, $false, $false, $false, $false, $true, 1, $false, "824÷6=137, 2", 2) | Out-Null
$d.Content.Find.Execute("731÷6=121, 5", $true, $false, $false, $false, $false, $true, 1, $false, "225÷6=37, 3", 2) | Out-Null
$d.Content.Find.Execute("325÷4=81, 1", $true, $false, $false, $false, $false, $true, 1, $false, "928÷4=232, 0", 2) | Out-Null
$d.Content.Find.Execute("586÷7=83, 5", $true, $false, $false, $false, $false, $true, 1, $false, "197÷4=49, 1", 2) | Out-Null
$d.Content.Find.Execute("910÷8=113, 6", $true, $false, $false, $false, $false, $true, 1, $false, "711÷4=177, 3", 2) | Out-Null
$d.Content.Find.Execute("811÷9=90, 1", $true, $false, $false, $false, $false, $true, 1, $false, "568÷4=142, 0", 2) | Out-Null
